$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1 -and $parts[0] -ceq "System") {
            $sortedParts = $parts | Sort-Object
            $newVal = $sortedParts -join ", "
            if ($newVal -cne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
